# edit.ps1
# Adds a new "2022-Q4" worksheet (fund holdings detail) positioned before the
# existing "2022-Q3" sheet, and updates the "总计" (totals) summary sheet with
# a new row for 2022-Q4, shifting the existing 2022-Q3 / 2022-Q2 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
#    Before: row2 = 2022-Q3 (66 / 16.83), row3 = 2022-Q2 (19 / 4.37)
#    After:  row2 = 2022-Q4 (21 / 8.01), row3 = 2022-Q3 (66 / 16.83),
#            row4 = 2022-Q2 (19 / 4.37)
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Make room for the new 2022-Q3 row: push the existing 2022-Q2 row (row 3)
# down to row 4.
$totals.Rows(3).Insert()

# Row 3 becomes the (old row 2) 2022-Q3 data.
$totals.Range("A3").Value = 1
$totals.Range("A3").Font.Bold = $true
$totals.Range("A3").HorizontalAlignment = -4108
$totals.Range("A3").VerticalAlignment = -4160
$totals.Range("A3").Borders.LineStyle = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 66
$totals.Range("D3").Value = 16.83

# Row 2 becomes the new 2022-Q4 summary data.
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 21
$totals.Range("D2").Value = 8.01

# Row 4 (shifted 2022-Q2 row) keeps its data but the index needs bumping
# from 1 to 2.
$totals.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fund holdings detail rows: code, name, size, position%, ratio%, marketValue, rank
$rows = @(
  @("180031","银华中小盘精选混合","37.32","91.61","4.62","1.7242",7),
  @("001822","华商智能生活灵活配置混合A","33.45","90.70","3.97","1.3280",8),
  @("007490","南方信息创新混合A","15.05","91.51","5.02","0.7555",8),
  @("506005","博时科创板三年定开混合","20.34","98.81","3.51","0.7139",8),
  @("001404","招商移动互联网产业股票A","13.83","90.58","4.78","0.6611",5),
  @("015385","华商智能生活灵活配置混合C","11.97","90.70","3.97","0.4752",8),
  @("013840","银华集成电路混合A","9.27","94.88","4.55","0.4218",7),
  @("015773","招商移动互联网产业股票C","8.53","90.58","4.78","0.4077",5),
  @("013841","银华集成电路混合C","8.03","94.88","4.55","0.3654",7),
  @("012650","博时半导体主题混合A","7.14","93.53","3.35","0.2392",10),
  @("162207","泰达宏利效率优选混合（LOF）","4.72","68.20","4.04","0.1907",7),
  @("009085","银华丰享一年持有期混合","3.60","91.96","4.60","0.1656",7),
  @("006864","国联安核心资产策略混合","4.67","91.48","3.45","0.1611",10),
  @("012651","博时半导体主题混合C","4.61","93.53","3.35","0.1544",10),
  @("007491","南方信息创新混合C","2.40","91.51","5.02","0.1205",8),
  @("009141","泰达宏利价值长青混合A","1.62","88.10","4.46","0.0723",8),
  @("015097","东财数字经济优选混合C","0.38","92.05","4.93","0.0187",6),
  @("015096","东财数字经济优选混合A","0.30","92.05","4.93","0.0148",6),
  @("015641","银华数字经济股票A","0.20","94.75","5.20","0.0104",7),
  @("009142","泰达宏利价值长青混合C","0.08","88.10","4.46","0.0036",8),
  @("015642","银华数字经济股票C","0.06","94.75","5.20","0.0031",7)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = 2 + $i
    $values = $rows[$i]

    $q4.Cells.Item($row, 1).Value = $i

    $codeCell = $q4.Cells.Item($row, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $values[0]

    $q4.Cells.Item($row, 3).Value = $values[1]

    $sizeCell = $q4.Cells.Item($row, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $values[2]

    $posCell = $q4.Cells.Item($row, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $values[3]

    $ratioCell = $q4.Cells.Item($row, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $values[4]

    $mvCell = $q4.Cells.Item($row, 7)
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $values[5]

    $q4.Cells.Item($row, 8).Value = $values[6]
}

$indexRange = $q4.Range("A2:A22")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1
